# Session 7: Branch and Bound
# Enter this student's Session 7 score and grader comment, which ripples
# through the weighted "Mark" (J4) and "Final mark" (K4) formulas already
# on the sheet, and leave behind the feedback note shown under the
# Session 7 header (row 5).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Session 7 (Branch&Bound) score for the student
$ws.Range("H4").Value = 9.5

# Per-session comment cell next to the score
$ws.Range("I4").Value = "na"

# Feedback text for Session 7, shown under the header row
$ws.Range("H5").Value = "You did a great work during the whole course. To greatly improve times, you may also comment lines 35 and 37 in Heap.java file. That part is only used to check if a node is already repeated but in this problem that will never happen. The times are going to be much much better without it. Your way of calculating the heuristic value is not so good: whenever you find a solution you will prune all the nodes which score lower than the solution you found (and those nodes may lead to a better solution after processing them)."

# Recalculate so the dependent "Mark"/"Final mark" formulas pick up the
# new Session 7 score.
$excel.Calculate()

# Move the view over to column H now that it is filled in, matching the
# way Excel scrolls the grid when you select the now-visible H5:H12 block.
$excel.Goto($ws.Range("H5:H12"), $true)
$excel.ActiveWindow.ScrollColumn = 2
